$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine last used row (xlUp)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Columns that may contain HYPERLINK formulas
$hyperlinkCols = @(19, 20, 21, 22, 23, 24, 25)  # S, T, U, V, W, X, Y

for ($r = 2; $r -le $lastRow; $r++) {
    # Update "Förändrad" date column (C = col 3)
    $cCell = $ws.Cells.Item($r, 3)
    if ($cCell.Value2 -eq 45184) {
        $cCell.Value2 = 45186
    }

    # Beteckning value used as the friendly name of hyperlinks
    $beteckning = $ws.Cells.Item($r, 1).Value2

    foreach ($col in $hyperlinkCols) {
        $cell = $ws.Cells.Item($r, $col)
        $formula = $cell.Formula
        if ($formula -and $formula.StartsWith("=HYPERLINK(")) {
            # Only add the friendly-name argument if it isn't already present
            if ($formula -notmatch ',\s*"') {
                $newFormula = $formula.Substring(0, $formula.Length - 1) + ', "' + $beteckning + '")'
                $cell.Formula = $newFormula
            }
        }
    }
}
